# Auto-generated Excel COM-interop script to apply the crypto price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) store text-formatted values (e.g. "512.36", "  +2.17%  ").
# Temporarily force Text number format so Excel does not reinterpret these as numeric values
# (which would round/alter the exact displayed text), then restore the original (unstyled) look.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Update Price (D) and Volume(1h) (E) columns for most rows ---
$ws.Range("D2").Value = "57.134.53"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "3.062.25"
$ws.Range("E3").Value = "  +4.74%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "512.36"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").Value = "142.25"
$ws.Range("E6").Value = "  +7.34%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +2.96%  "
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("E10").Value = "  +3.54%  "
$ws.Range("E11").Value = "  +5.88%  "
$ws.Range("D12").Value = "3.592.44"
$ws.Range("E12").Value = "  +5.08%  "
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").Value = "25.63"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("E15").Value = "  +2.88%  "
$ws.Range("D16").Value = "57.301.61"
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").Value = "3.056.67"
$ws.Range("E17").Value = "  +4.82%  "
$ws.Range("D18").Value = "6.08"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").Value = "12.94"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").Value = "8.15"
$ws.Range("E20").Value = "  +5.91%  "
$ws.Range("D21").Value = "335.15"
$ws.Range("E21").Value = "  +6.77%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "0.497"
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("D24").Value = "65.24"
$ws.Range("E24").Value = "  +3.89%  "
$ws.Range("E25").Value = "  +7.22%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").Value = "0.0₃0924"
$ws.Range("E27").Value = "  +10.28%  "
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("D29").Value = "7.04"
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").Value = "20.69"
$ws.Range("E31").Value = "  +4.09%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").Value = "154.18"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("D34").Value = "4.50"
$ws.Range("E34").Value = "  +3.14%  "
$ws.Range("D35").Value = "5.87"
$ws.Range("E35").Value = "  +5.31%  "
$ws.Range("D36").Value = "26.23"
$ws.Range("E36").Value = "  +9.61%  "
$ws.Range("E37").Value = "  +4.02%  "
$ws.Range("D38").Value = "0.0680"
$ws.Range("E38").Value = "  +5.67%  "
$ws.Range("D39").Value = "3.102.63"
$ws.Range("E39").Value = "  +4.95%  "
$ws.Range("D40").Value = "36.65"
$ws.Range("E40").Value = "  +0.69%  "

# --- Rows 41-43: coin ordering changed (Filecoin / FirstDigitalUSD / Mantle reshuffled) ---
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.667"
$ws.Range("E42").Value = "  +4.73%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.83"
$ws.Range("E43").Value = "  +3.53%  "

$ws.Range("D44").Value = "2.263.76"
$ws.Range("E44").Value = "  +7.07%  "
$ws.Range("E45").Value = "  +8.50%  "
$ws.Range("E46").Value = "  +3.88%  "
$ws.Range("E47").Value = "  +3.50%  "
$ws.Range("D48").Value = "20.25"
$ws.Range("E48").Value = "  +8.63%  "
$ws.Range("D49").Value = "5.84"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("E50").Value = "  +4.54%  "
$ws.Range("E51").Value = "  +7.24%  "

# Restore the original (default/unstyled) appearance of the Price/Volume columns
# now that the text values are safely in place.
$ws.Range("D2:E51").Style = "Normal"

